# Insert the new explanatory sentence about creating the ingress namespace
# into the last (empty) paragraph of the document body, right before the
# final section break — mirrors the target diff, which adds a run of text
# (interspersed with Word's auto spell-check <w:proofErr/> markers around
# the technical terms "ingress", "namespace", "ingress-nginx" and
# "kubectl") into that trailing empty paragraph.

$d = $word.ActiveDocument

# The paragraph to edit is the very last paragraph in the main document
# story (an empty paragraph that sits immediately before the final
# sectPr/section break).
$target = $d.Paragraphs.Last.Range

$newText = "Para hacer el ingress , tendremos que crear el namespace de tipo ingress-nginx con kubectl, luego "

$target.InsertBefore($newText)
